$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Pre-format column B as Text so the long hash bit-strings are not
# coerced into numbers (scientific notation / precision loss).
$ws.Range("B1:B63").NumberFormat = "@"

# Header row
$ws.Range("B1").Value = "Хэш"
$ws.Range("C1").Value = "Время обработки"
$ws.Range("D1").Value = "Хэммингово расстояние"

$ws.Range("B2").Value = "1111111111011011100010010000000100001101100011011101101111111111"
$ws.Range("C2").Value = 0.020635
$ws.Range("D2").Value = 0

$ws.Range("B3").Value = "1111111111011011100010010000000100001101100011111101101111111111"
$ws.Range("C3").Value = 0.0255
$ws.Range("D3").Value = 1

$ws.Range("B4").Value = "1111111111011011100010010000000100001101100011011101101111111111"
$ws.Range("C4").Value = 0.007990000000000001
$ws.Range("D4").Value = 0

$ws.Range("B5").Value = "1111111111011011100010010000000100001101100011011101101111111111"
$ws.Range("C5").Value = 0.032334
$ws.Range("D5").Value = 0

$ws.Range("B6").Value = "1111111111011011100010010000000100001101100011111101101111111111"
$ws.Range("C6").Value = 0
$ws.Range("D6").Value = 1

$ws.Range("B7").Value = "1111111111011011100010010000000100001101100011111101101111111111"
$ws.Range("C7").Value = 0
$ws.Range("D7").Value = 1

$ws.Range("B8").Value = "1111111111011011100010010000000100001101100011111101101111111111"
$ws.Range("C8").Value = 0
$ws.Range("D8").Value = 1

$ws.Range("B9").Value = "1111111111011011100010010000000100001101100011111101101111111111"
$ws.Range("C9").Value = 0
$ws.Range("D9").Value = 1

$ws.Range("B10").Value = "1111111111011011100010010000000100001101100011111101101111111111"
$ws.Range("C10").Value = 0
$ws.Range("D10").Value = 1

$ws.Range("B11").Value = "1111111111011011100010010000000100001101100011111101101111111111"
$ws.Range("C11").Value = 0
$ws.Range("D11").Value = 1

$ws.Range("B12").Value = "1111111111011011100010010000000100001101100011011101101111111111"
$ws.Range("C12").Value = 0
$ws.Range("D12").Value = 0

$ws.Range("B13").Value = "1111111111011011100010010000000100001101100011011101101111111111"
$ws.Range("C13").Value = 0
$ws.Range("D13").Value = 0

$ws.Range("B14").Value = "1111111111011011100010010000000100001101100011011101101111111111"
$ws.Range("C14").Value = 0
$ws.Range("D14").Value = 0

$ws.Range("B15").Value = "1111111111011011100010010000000100001101100011011101101111111111"
$ws.Range("C15").Value = 0
$ws.Range("D15").Value = 0

$ws.Range("B16").Value = "1111111111011011100010010000000100001101100011011101101111111111"
$ws.Range("C16").Value = 0
$ws.Range("D16").Value = 0

$ws.Range("B17").Value = "1111111111011011100010010000000100001101100011011101101111111111"
$ws.Range("C17").Value = 0
$ws.Range("D17").Value = 0

$ws.Range("B18").Value = "1111111111011011100010010000000100001101100011011101101111111111"
$ws.Range("C18").Value = 0
$ws.Range("D18").Value = 0

$ws.Range("B19").Value = "1111111111011011100010010000000100001101100011011101101111111111"
$ws.Range("C19").Value = 0
$ws.Range("D19").Value = 0

$ws.Range("B20").Value = "1111111111011011100010010000000100001101100011011101101111111111"
$ws.Range("C20").Value = 0.015617
$ws.Range("D20").Value = 0

$ws.Range("B21").Value = "1111111111011011100010010000000100001101100011011101101111111111"
$ws.Range("C21").Value = 0
$ws.Range("D21").Value = 0

$ws.Range("B22").Value = "1111111111011011100010010000000100001101100011011101101111111111"
$ws.Range("C22").Value = 0
$ws.Range("D22").Value = 0

$ws.Range("B23").Value = "1111111111011011100010010000000100001101100011011101101111111111"
$ws.Range("C23").Value = 0.015637
$ws.Range("D23").Value = 0

$ws.Range("B24").Value = "1111111111011011100010010000000100001101100011011101101111111111"
$ws.Range("C24").Value = 0.014129
$ws.Range("D24").Value = 0

$ws.Range("B25").Value = "1111111111011011100010010000000100001101100011111101111111111111"
$ws.Range("C25").Value = 0.007818
$ws.Range("D25").Value = 2

$ws.Range("B26").Value = "1111111111011011100011110000001100001101100011111101111111111111"
$ws.Range("C26").Value = 0.015693
$ws.Range("D26").Value = 5

$ws.Range("B27").Value = "1111111111011111100011110000101100001111100011111101111111111111"
$ws.Range("C27").Value = 0.015622
$ws.Range("D27").Value = 8

$ws.Range("B28").Value = "1111111111011111100011110000111100001111100011111101111111111111"
$ws.Range("C28").Value = 0.015625
$ws.Range("D28").Value = 9

$ws.Range("B29").Value = "1111111111011111100011110000111100001111100011111101111111111111"
$ws.Range("C29").Value = 0.015632
$ws.Range("D29").Value = 9

$ws.Range("B30").Value = "1111111111011111100011110000111100101111100011111101111111111111"
$ws.Range("C30").Value = 0.015621
$ws.Range("D30").Value = 10

$ws.Range("B31").Value = "1111111111011111100011110000111100101111000011111101111111111111"
$ws.Range("C31").Value = 0.015649
$ws.Range("D31").Value = 11

$ws.Range("B32").Value = "1111111111011111100011110000111100101111000011111101111111111111"
$ws.Range("C32").Value = 0
$ws.Range("D32").Value = 11

$ws.Range("B33").Value = "1111111111011111000011110000111100101111000011111101111111111111"
$ws.Range("C33").Value = 0.015612
$ws.Range("D33").Value = 12

$ws.Range("B34").Value = "1111111111100111110001111100001111000111110001111110011111111111"
$ws.Range("C34").Value = 0.015615
$ws.Range("D34").Value = 22

$ws.Range("B35").Value = "1110001111001001101000011110000111101101100010011100100111100011"
$ws.Range("C35").Value = 0
$ws.Range("D35").Value = 19

$ws.Range("B36").Value = "1111011111000011100000011110000110000001100000011100001111110111"
$ws.Range("C36").Value = 0
$ws.Range("D36").Value = 15

$ws.Range("B37").Value = "1111111111010011100100011001000011110000111110011111001111110111"
$ws.Range("C37").Value = 0.016628
$ws.Range("D37").Value = 20

$ws.Range("B38").Value = "1111101110000001101100001011000010011001100110011101101111111111"
$ws.Range("C38").Value = 0
$ws.Range("D38").Value = 18

$ws.Range("B39").Value = "1111111111111011100000111000011110000111100001111100111111111111"
$ws.Range("C39").Value = 0
$ws.Range("D39").Value = 13

$ws.Range("B40").Value = "1111111111000011100000011000000110000001100000011100001111111111"
$ws.Range("C40").Value = 0.015626
$ws.Range("D40").Value = 11

$ws.Range("B41").Value = "1111111111000101100001011000111110001111100011111101111111111111"
$ws.Range("C41").Value = 0
$ws.Range("D41").Value = 14

$ws.Range("B42").Value = "1111111111111111110001111100001110000011110000111110011111111111"
$ws.Range("C42").Value = 0.015627
$ws.Range("D42").Value = 21

$ws.Range("B43").Value = "1111111111000001100000011000100110001101100010011100001111111111"
$ws.Range("C43").Value = 0
$ws.Range("D43").Value = 10

$ws.Range("B44").Value = "1111111111000011100000111000001111001011111100111110011111111111"
$ws.Range("C44").Value = 0
$ws.Range("D44").Value = 20

$ws.Range("B45").Value = "1111111111000001100000011000010110001101100011011100001111111111"
$ws.Range("C45").Value = 0.015623
$ws.Range("D45").Value = 9

$ws.Range("B46").Value = "1111011111101011110000111000000111001011110010111110101111100111"
$ws.Range("C46").Value = 0
$ws.Range("D46").Value = 18

$ws.Range("B47").Value = "1111111111001011100001111000001110000001100111111111111111111111"
$ws.Range("C47").Value = 0
$ws.Range("D47").Value = 13

$ws.Range("B48").Value = "1111011111000011100011011010010100000101100011111100111111101111"
$ws.Range("C48").Value = 0
$ws.Range("D48").Value = 12

$ws.Range("B49").Value = "1110011111000011110010011100100111001001110010011100001111100011"
$ws.Range("C49").Value = 0
$ws.Range("D49").Value = 18

$ws.Range("B50").Value = "1111000111111001110000011100010111000001110000011100101111000011"
$ws.Range("C50").Value = 0.015624
$ws.Range("D50").Value = 22

$ws.Range("B51").Value = "1111111111001011100000010000010100000001100011111101111111111111"
$ws.Range("C51").Value = 0
$ws.Range("D51").Value = 7

$ws.Range("B52").Value = "1111111111001011100000010000010100000001100011111101111111111111"
$ws.Range("C52").Value = 0.015837
$ws.Range("D52").Value = 7

$ws.Range("B53").Value = "1111111111001011100001011000111110001111100011111111111111111111"
$ws.Range("C53").Value = 0
$ws.Range("D53").Value = 12

$ws.Range("B54").Value = "1111101111001001100001010010000000000101000011011100110111101011"
$ws.Range("C54").Value = 0.0157
$ws.Range("D54").Value = 14

$ws.Range("B55").Value = "1111111111110011110000111000000110000001110000111110001111111111"
$ws.Range("C55").Value = 0.015622
$ws.Range("D55").Value = 16

$ws.Range("B56").Value = "1101111110011111001100001011000000110000101100001111101111111111"
$ws.Range("C56").Value = 0.015581
$ws.Range("D56").Value = 23

$ws.Range("B57").Value = "1110011111000011110000111100001111000011110000111100011111100111"
$ws.Range("C57").Value = 0
$ws.Range("D57").Value = 24

$ws.Range("B58").Value = "1111111111110011100000011000000100000001100010011111101111111111"
$ws.Range("C58").Value = 0.01562
$ws.Range("D58").Value = 8

$ws.Range("B59").Value = "1111111111000111110000111100001111000011110000111100011111101111"
$ws.Range("C59").Value = 0.031308
$ws.Range("D59").Value = 22

$ws.Range("B60").Value = "1111101111000011100001111000000110000001111000111111001111111111"
$ws.Range("C60").Value = 0.031244
$ws.Range("D60").Value = 17

$ws.Range("B61").Value = "1110001111000001100000000000010000000100100011011100110111100011"
$ws.Range("C61").Value = 0.03125
$ws.Range("D61").Value = 18

$ws.Range("B62").Value = "1111101111001011100001111000001100000011100000111000011111111111"
$ws.Range("C62").Value = 0.031201
$ws.Range("D62").Value = 17

$ws.Range("B63").Value = "1111111111110111110000111100001111000011110000111111011111111111"
$ws.Range("C63").Value = 0.046921
$ws.Range("D63").Value = 21
